# Add "purging of channel before stream" credentials columns (API_KEY, CHANNEL_ID)
# to the Credentials sheet, and update the active sheet/selection state to match.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("test_schedule_1")
$wsCreds    = $wb.Worksheets.Item("Credentials")

# --- Credentials sheet: new columns E (API_KEY) and F (CHANNEL_ID) ---
$wsCreds.Range("E1").Value = "API_KEY"
$wsCreds.Range("F1").Value = "CHANNEL_ID"
$wsCreds.Range("E2").Value = 12345
$wsCreds.Range("F2").Value = 12345

# Give the new CHANNEL_ID column the same custom width treatment the sheet
# already uses for its other data columns (closest width reachable through
# the ColumnWidth property's character-width rounding).
$wsCreds.Columns.Item(6).ColumnWidth = 12

# --- Selections on each sheet ---
$wsSchedule.Range("C10").Select() | Out-Null
$wsCreds.Range("L7").Select() | Out-Null

# --- Make the Credentials sheet the active/selected tab ---
$wsCreds.Activate() | Out-Null
